# Updates cryptos list price (D) and volume/change (E) columns
# per the commit "Updated cryptos list on Wed Aug 21 12:51:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.429.01"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "'2.574.77"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'552.03"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "'2.580.91"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").Value = "  +8.70%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "'3.028.66"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "'59.412.26"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "'23.10"
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'2.579.62"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'338.39"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'10.32"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'0.481"
$ws.Range("E24").Value = "  +7.61%  "
$ws.Range("E25").Value = "  -5.22%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").Value = "'7.39"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'0.0₃0772"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'6.19"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "'158.79"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").Value = "'19.04"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").Value = "'0.901"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").Value = "'37.45"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'0.852"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "'289.30"
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").Value = "'138.24"
$ws.Range("E43").Value = "  +8.26%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.0970"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'0.592"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'0.0530"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'18.65"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'1.952.96"
$ws.Range("E51").Value = "  -0.53%  "
